$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1558
$ws.Range("I41").Value = 1400
$ws.Range("J41").Value = 1663.3334
$ws.Range("K41").Value = 1400
$ws.Range("L41").Value = 1663.3334
$ws.Range("M41").Value = -960
$ws.Range("N41").Value = -2543.3334

$ws.Range("H63").Value = 39447.332
$ws.Range("J63").Value = 39447.332
$ws.Range("L63").Value = 39447.332
$ws.Range("N63").Value = -40695.332

$ws.Range("H66").Value = 39447.332
$ws.Range("J66").Value = 39447.332
$ws.Range("L66").Value = 118341.996
$ws.Range("N66").Value = -124581.996

$ws.Range("H74").Value = 4468.2856
$ws.Range("I74").Value = 4172.727
$ws.Range("J74").Value = 5552
$ws.Range("K74").Value = 4172.727
$ws.Range("L74").Value = 5552
$ws.Range("M74").Value = -3236.727
$ws.Range("N74").Value = -7424

$ws.Range("H77").Value = 4468.2856
$ws.Range("I77").Value = 4172.727
$ws.Range("J77").Value = 5552
$ws.Range("K77").Value = 20863.635
$ws.Range("L77").Value = 27760
$ws.Range("M77").Value = -16183.635
$ws.Range("N77").Value = -37120

$ws.Range("H112").Value = 1908.7059
$ws.Range("I112").Value = 527.5
$ws.Range("J112").Value = 2333.6924
$ws.Range("K112").Value = 1582.5
$ws.Range("L112").Value = 7001.0772
$ws.Range("M112").Value = -474.5
$ws.Range("N112").Value = -9217.0772

$ws.Range("H137").Value = 1258.7646
$ws.Range("I137").Value = 1213.3529
$ws.Range("J137").Value = 1349.5883
$ws.Range("K137").Value = 3640.0587
$ws.Range("L137").Value = 4048.7649
$ws.Range("M137").Value = -1090.0587
$ws.Range("N137").Value = -9148.7649

$ws.Range("H138").Value = 3873.0454
$ws.Range("I138").Value = 823.5833
$ws.Range("J138").Value = 7532.4
$ws.Range("K138").Value = 2470.7499
$ws.Range("L138").Value = 22597.2
$ws.Range("M138").Value = 2669.2501
$ws.Range("N138").Value = -32877.2

$ws.Range("H141").Value = 1718.7826
$ws.Range("I141").Value = 1706
$ws.Range("J141").Value = 2000
$ws.Range("K141").Value = 5118
$ws.Range("L141").Value = 6000
$ws.Range("M141").Value = 62
$ws.Range("N141").Value = -16360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3803.8408
$ws.Range("I61").Value = 4037.1843
$ws.Range("J61").Value = 2326
$ws.Range("K61").Value = 4037.1843
$ws.Range("L61").Value = 2326
$ws.Range("M61").Value = -3825.1843
$ws.Range("N61").Value = -2750

$ws.Range("H122").Value = 1976006.9
$ws.Range("I122").Value = 2853653
$ws.Range("J122").Value = 1303
$ws.Range("K122").Value = 8560959
$ws.Range("L122").Value = 3909
$ws.Range("M122").Value = -8558509
$ws.Range("N122").Value = -8809

$ws.Range("H132").Value = 3152.5527
$ws.Range("I132").Value = 1784.4615
$ws.Range("J132").Value = 6116.75
$ws.Range("K132").Value = 5353.3845
$ws.Range("L132").Value = 18350.25
$ws.Range("M132").Value = -2823.3845
$ws.Range("N132").Value = -23410.25

$ws.Range("H136").Value = 3803.8408
$ws.Range("I136").Value = 4037.1843
$ws.Range("J136").Value = 2326
$ws.Range("K136").Value = 12111.5529
$ws.Range("L136").Value = 6978
$ws.Range("M136").Value = -9561.552899999999
$ws.Range("N136").Value = -12078

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H93").Value = 60000
$ws.Range("J93").Value = 60000
$ws.Range("L93").Value = 60000
$ws.Range("N93").Value = -63744

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3624426.5
$ws.Range("I58").Value = 7247300.5
$ws.Range("J58").Value = 1552.3043
$ws.Range("K58").Value = 7247300.5
$ws.Range("L58").Value = 1552.3043
$ws.Range("M58").Value = -7247097.5
$ws.Range("N58").Value = -1958.3043

$ws.Range("H62").Value = 5440
$ws.Range("I62").Value = 5562.5
$ws.Range("J62").Value = 4950
$ws.Range("K62").Value = 5562.5
$ws.Range("L62").Value = 4950
$ws.Range("M62").Value = -4938.5
$ws.Range("N62").Value = -6198

$ws.Range("H65").Value = 5440
$ws.Range("I65").Value = 5562.5
$ws.Range("J65").Value = 4950
$ws.Range("K65").Value = 27812.5
$ws.Range("L65").Value = 24750
$ws.Range("M65").Value = -24692.5
$ws.Range("N65").Value = -30990

$ws.Range("H87").Value = 21975
$ws.Range("J87").Value = 21975
$ws.Range("L87").Value = 21975
$ws.Range("N87").Value = -24347

$ws.Range("H90").Value = 21975
$ws.Range("J90").Value = 21975
$ws.Range("L90").Value = 65925
$ws.Range("N90").Value = -77781

$ws.Range("H99").Value = 10422632
$ws.Range("I99").Value = 4158.8
$ws.Range("K99").Value = 4158.8
$ws.Range("M99").Value = -2660.8

$ws.Range("H126").Value = 10422632
$ws.Range("I126").Value = 4158.8
$ws.Range("K126").Value = 12476.4
$ws.Range("M126").Value = -10006.4

$ws.Range("H136").Value = 3624426.5
$ws.Range("I136").Value = 7247300.5
$ws.Range("J136").Value = 1552.3043
$ws.Range("K136").Value = 21741901.5
$ws.Range("L136").Value = 4656.9129
$ws.Range("M136").Value = -21739351.5
$ws.Range("N136").Value = -9756.912899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 4846129
$ws.Range("J5").Value = 10768254
$ws.Range("L5").Value = 32304762
$ws.Range("N5").Value = -32304986

$ws.Range("H135").Value = 4846129
$ws.Range("J135").Value = 10768254
$ws.Range("L135").Value = 96914286
$ws.Range("N135").Value = -96919356

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6178.0713
$ws.Range("I70").Value = 6162.591
$ws.Range("J70").Value = 6234.8335
$ws.Range("K70").Value = 6162.591
$ws.Range("L70").Value = 6234.8335
$ws.Range("M70").Value = -5892.591
$ws.Range("N70").Value = -6774.8335

$ws.Range("H73").Value = 6178.0713
$ws.Range("I73").Value = 6162.591
$ws.Range("J73").Value = 6234.8335
$ws.Range("K73").Value = 6162.591
$ws.Range("L73").Value = 6234.8335
$ws.Range("M73").Value = -5226.591
$ws.Range("N73").Value = -8106.8335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2175.375
$ws.Range("I7").Value = 1980.6
$ws.Range("K7").Value = 1980.6
$ws.Range("M7").Value = -1868.6

$ws.Range("H40").Value = 27779000
$ws.Range("I40").Value = 50001140
$ws.Range("J40").Value = 1326.5625
$ws.Range("K40").Value = 50001140
$ws.Range("L40").Value = 1326.5625
$ws.Range("M40").Value = -50001004
$ws.Range("N40").Value = -1598.5625

$ws.Range("H126").Value = 2175.375
$ws.Range("I126").Value = 1980.6
$ws.Range("K126").Value = 5941.799999999999
$ws.Range("M126").Value = -3471.799999999999

$ws.Range("H133").Value = 40320
$ws.Range("J133").Value = 40320
$ws.Range("L133").Value = 40320
$ws.Range("N133").Value = -45380

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 50000
$ws.Range("J82").Value = 50000
$ws.Range("L82").Value = 50000
$ws.Range("N82").Value = -50766

$ws.Range("H85").Value = 50000
$ws.Range("J85").Value = 50000
$ws.Range("L85").Value = 50000
$ws.Range("N85").Value = -52652

$ws.Range("H96").Value = 2605.9092
$ws.Range("I96").Value = 2270
$ws.Range("J96").Value = 3193.75
$ws.Range("K96").Value = 2270
$ws.Range("L96").Value = 3193.75
$ws.Range("M96").Value = -897
$ws.Range("N96").Value = -5939.75

$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws.Range("H122").Value = 998.3333
$ws.Range("I122").Value = 999
$ws.Range("J122").Value = 997.5
$ws.Range("K122").Value = 2997
$ws.Range("L122").Value = 2992.5
$ws.Range("M122").Value = -547
$ws.Range("N122").Value = -7892.5

$ws.Range("H126").Value = 875.1875
$ws.Range("I126").Value = 666.0833
$ws.Range("J126").Value = 1502.5
$ws.Range("K126").Value = 1998.2499
$ws.Range("L126").Value = 4507.5
$ws.Range("M126").Value = 471.7501
$ws.Range("N126").Value = -9447.5

$ws.Range("H132").Value = 1652.6538
$ws.Range("I132").Value = 1376.0625
$ws.Range("J132").Value = 2095.2
$ws.Range("K132").Value = 4128.1875
$ws.Range("L132").Value = 6285.599999999999
$ws.Range("M132").Value = -1598.1875
$ws.Range("N132").Value = -11345.6
